$d = $word.ActiveDocument

# --- Paragraph 158 ("Throughout the project, ... thoroughly documented") ---
# Fix typo "member" -> "members" in "neither new member nor"
$p = $d.Paragraphs(158)
$p.Range.Find.Execute("neither new member nor", $true, $false, $false, $false, $false, $true, 1, $false, "neither new members nor", 1) | Out-Null

# Split off the trailing (empty) run into its own paragraph
$p = $d.Paragraphs(158)
$p.Range.Find.Execute("thoroughly documented", $true, $false, $false, $false, $false, $true, 1, $false, "thoroughly documented^p", 1) | Out-Null

# Re-indent both the (now shorter) original paragraph and the new split paragraph
$p = $d.Paragraphs(158)
$p.Format.LeftIndent = 36
$p.Format.FirstLineIndent = 36
$pNew = $d.Paragraphs(159)
$pNew.Format.LeftIndent = 36
$pNew.Format.FirstLineIndent = 0

Write-Output "done"
